$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.924.72"
$ws.Range("E2").Value = "  +0.24%  "

Set-TextValue $ws.Range("D3") "2.357.95"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "0.692"
$ws.Range("E5").Value = "  +4.53%  "

Set-TextValue $ws.Range("D6") "241.97"
$ws.Range("E6").Value = "  +2.81%  "

Set-TextValue $ws.Range("D7") "77.26"
$ws.Range("E7").Value = "  +5.19%  "

$ws.Range("E8").Value = "  -0.01%  "

Set-TextValue $ws.Range("D9") "0.633"
$ws.Range("E9").Value = "  +21.36%  "

$ws.Range("E10").Value = "  +3.56%  "

Set-TextValue $ws.Range("D11") "57.40"
$ws.Range("E11").Value = "  +0.74%  "

Set-TextValue $ws.Range("D12") "33.65"
$ws.Range("E12").Value = "  +22.86%  "

Set-TextValue $ws.Range("D13") "7.55"
$ws.Range("E13").Value = "  +14.65%  "

$ws.Range("E14").Value = "  +1.72%  "

Set-TextValue $ws.Range("D15") "2.708.24"
$ws.Range("E15").Value = "  -0.45%  "

Set-TextValue $ws.Range("D16") "17.03"
$ws.Range("E16").Value = "  +3.02%  "

Set-TextValue $ws.Range("D17") "0.934"
$ws.Range("E17").Value = "  +5.09%  "

Set-TextValue $ws.Range("D18") "2.358.30"
$ws.Range("E18").Value = "  -0.52%  "

Set-TextValue $ws.Range("D19") "43.810.76"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("E20").Value = "  +1.84%  "

Set-TextValue $ws.Range("D21") "6.69"
$ws.Range("E21").Value = "  +3.51%  "

Set-TextValue $ws.Range("D22") "77.75"
$ws.Range("E22").Value = "  +2.53%  "

Set-TextValue $ws.Range("D23") "259.63"
$ws.Range("E23").Value = "  +3.06%  "

$ws.Range("E24").Value = "  +0.06%  "

Set-TextValue $ws.Range("D25") "2.54"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("E26").Value = "  -3.78%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "11.03"
$ws.Range("E27").Value = "  +7.29%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D28") "1.81"
$ws.Range("E28").Value = "  +17.76%  "

$ws.Range("E29").Value = "  +2.31%  "

Set-TextValue $ws.Range("D30") "23.15"
$ws.Range("E30").Value = "  +2.25%  "

Set-TextValue $ws.Range("D31") "175.11"
$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("E32").Value = "  -4.53%  "

$ws.Range("E33").Value = "  +4.31%  "

Set-TextValue $ws.Range("D34") "0.0768"
$ws.Range("E34").Value = "  +9.27%  "

Set-TextValue $ws.Range("D35") "5.38"
$ws.Range("E35").Value = "  +4.95%  "

Set-TextValue $ws.Range("D36") "5.41"
$ws.Range("E36").Value = "  +4.99%  "

Set-TextValue $ws.Range("D37") "3.80"
$ws.Range("E37").Value = "  +0.52%  "

Set-TextValue $ws.Range("D38") "2.44"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("E39").Value = "  -3.21%  "

Set-TextValue $ws.Range("D40") "0.0280"
$ws.Range("E40").Value = "  +6.56%  "

Set-TextValue $ws.Range("D41") "0.110"
$ws.Range("E41").Value = "  +14.03%  "

Set-TextValue $ws.Range("D42") "0.208"
$ws.Range("E42").Value = "  +17.90%  "

Set-TextValue $ws.Range("D43") "19.40"
$ws.Range("E43").Value = "  -1.52%  "

Set-TextValue $ws.Range("D44") "9.14"
$ws.Range("E44").Value = "  +2.93%  "

$ws.Range("E45").Value = "  -0.15%  "

Set-TextValue $ws.Range("D46") "2.55"
$ws.Range("E46").Value = "  +11.59%  "

Set-TextValue $ws.Range("D47") "1.27"
$ws.Range("E47").Value = "  +4.01%  "

Set-TextValue $ws.Range("D48") "1.20"
$ws.Range("E48").Value = "  +1.88%  "

Set-TextValue $ws.Range("D49") "102.12"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D50") "4.54"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D51") "56.57"
$ws.Range("E51").Value = "  +10.61%  "
